$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.787|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.213|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B3").Value = "x1:0.000|x2:0.660|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.340|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B4").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:1.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B5").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.455|x5:0.545|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B6").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:1.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B7").Value = "x1:0.000|x2:0.164|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.294|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.359|x13:0.000|x14:0.183|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B8").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:1.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B9").Value = "x1:0.000|x2:0.000|x3:1.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B10").Value = "x1:0.000|x2:0.538|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.462|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B11").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.517|x12:0.483|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B12").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.312|x10:0.365|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.324|x17:0.000"
$ws.Range("B13").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:1.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B14").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:1.000|x17:0.000"
$ws.Range("B15").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:1.000|x16:0.000|x17:0.000"
$ws.Range("B16").Value = "x1:0.719|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.281|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B17").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:1.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B18").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.481|x11:0.000|x12:0.000|x13:0.519|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B19").Value = "x1:0.719|x2:1.362|x3:1.000|x4:0.455|x5:0.545|x6:1.000|x7:1.081|x8:1.000|x9:1.312|x10:1.467|x11:0.517|x12:0.842|x13:1.732|x14:1.645|x15:1.000|x16:1.324|x17:0.000"
